$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 32x2 data block (columns A and B) with new computed values
$ws.Range("A1").Value = -0.29884856315290165
$ws.Range("B1").Value = 0.29817188070433076
$ws.Range("A2").Value = -0.19289295398641215
$ws.Range("B2").Value = 0.19134160112085397
$ws.Range("A3").Value = -0.088395959210014752
$ws.Range("B3").Value = 0.088121052006904677
$ws.Range("A4").Value = -0.076121052090261543
$ws.Range("B4").Value = 0.075905969862031242
$ws.Range("A5").Value = -0.14429477392719736
$ws.Range("B5").Value = 0.14267708450299654
$ws.Range("A6").Value = -0.10941671732077118
$ws.Range("B6").Value = 0.10923475393768056
$ws.Range("A7").Value = -0.089234754326978916
$ws.Range("B7").Value = 0.088775749813247984
$ws.Range("A8").Value = -0.068775750207324293
$ws.Range("B8").Value = 0.068378874490827357
$ws.Range("A9").Value = -0.062378874827482278
$ws.Range("B9").Value = 0.062040904471547975
$ws.Range("A10").Value = -0.056040904813173142
$ws.Range("B10").Value = 0.05599213864628183
$ws.Range("A11").Value = -0.051492138981888758
$ws.Range("B11").Value = 0.051410300061899505
$ws.Range("A12").Value = -0.045410300405622106
$ws.Range("B12").Value = 0.045157161924044953
$ws.Range("A13").Value = -0.03915716227318633
$ws.Range("B13").Value = 0.03908821468193846
$ws.Range("A14").Value = -0.027088215059818843
$ws.Range("B14").Value = 0.027054849069315878
$ws.Range("A15").Value = -0.021054849421450861
$ws.Range("B15").Value = 0.021028586414063355
$ws.Range("A16").Value = -0.015028586767426466
$ws.Range("B16").Value = 0.015004623271638096
$ws.Range("A17").Value = -0.009004623626594821
$ws.Range("B17").Value = 0.0089999996311158625
$ws.Range("A18").Value = -0.036111826212465559
$ws.Range("B18").Value = 0.036096928185106236
$ws.Range("A19").Value = -0.027096928501912476
$ws.Range("B19").Value = 0.027013793592558955
$ws.Range("A20").Value = -0.018013793912317055
$ws.Range("B20").Value = 0.018004281202562922
$ws.Range("A21").Value = -0.0090042815227997508
$ws.Range("B21").Value = 0.0089999996794096759
$ws.Range("A22").Value = -0.09394564223508084
$ws.Range("B22").Value = 0.093633155160585346
$ws.Range("A23").Value = -0.084633155487926714
$ws.Range("B23").Value = 0.084126583769045027
$ws.Range("A24").Value = -0.042126584250611998
$ws.Range("B24").Value = 0.041999999515674347
$ws.Range("A25").Value = -0.027260367506841021
$ws.Range("B25").Value = 0.027228500980669423
$ws.Range("A26").Value = -0.074414272494124134
$ws.Range("B26").Value = 0.074220045505207111
$ws.Range("A27").Value = -0.068220045831815845
$ws.Range("B27").Value = 0.067572259629001152
$ws.Range("A28").Value = -0.061572259961524267
$ws.Range("B28").Value = 0.061136938332380986
$ws.Range("A29").Value = -0.049136938695898635
$ws.Range("B29").Value = 0.048938766312941695
$ws.Range("A30").Value = -0.028938766714436515
$ws.Range("B30").Value = 0.028882411442479494
$ws.Range("A31").Value = -0.059659786189937947
$ws.Range("B31").Value = 0.059531401270005446
$ws.Range("A32").Value = -0.038531401679893662
$ws.Range("B32").Value = 0.038459639069819751

# Narrow columns A and B by one character each (16.43 -> 15.43, 15.71 -> 14.71)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666668
$ws.Columns.Item(2).ColumnWidth = 13.833333333333332
